# Updated cryptos list on Mon Aug 19 16:01:00 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row.
# Several Price values are decimal strings that Excel's type-inference would
# otherwise coerce into numbers (dropping significant trailing zeros, e.g.
# "543.50" -> 543.5), so those are written with a leading apostrophe to force
# them to stay literal text, exactly like typing them into Excel would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '58.500.74'
$ws.Range("E2").Value = '  -2.14%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '2.580.57'
$ws.Range("E3").Value = '  -3.40%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.09%  '

# Row 5: BNB
$ws.Range("D5").Value = '''543.50'

# Row 6: Solana
$ws.Range("D6").Value = '''144.00'
$ws.Range("E6").Value = '  -1.19%  '

# Row 7: USDC
$ws.Range("E7").Value = '  -0.03%  '

# Row 8: XRP
$ws.Range("E8").Value = '  +1.06%  '

# Row 9: Toncoin
$ws.Range("D9").Value = '''6.75'
$ws.Range("E9").Value = '  +1.19%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  -3.16%  '

# Row 11: TRON
$ws.Range("E11").Value = '  +3.53%  '

# Row 12: Cardano
$ws.Range("E12").Value = '  -1.73%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '3.037.61'
$ws.Range("E13").Value = '  -3.04%  '

# Row 14: WrappedBTC
$ws.Range("D14").Value = '58.449.35'
$ws.Range("E14").Value = '  -2.11%  '

# Row 15: Avalanche
$ws.Range("D15").Value = '''20.61'
$ws.Range("E15").Value = '  -2.89%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '2.596.21'
$ws.Range("E16").Value = '  -1.75%  '

# Row 17: ShibaInu
$ws.Range("D17").Value = '''0.0000131'
$ws.Range("E17").Value = '  -3.01%  '

# Row 18: Polkadot
$ws.Range("D18").Value = '''4.45'

# Row 19: BitcoinCash
$ws.Range("D19").Value = '''334.46'
$ws.Range("E19").Value = '  -3.02%  '

# Row 20: Chainlink
$ws.Range("D20").Value = '''10.05'
$ws.Range("E20").Value = '  -3.19%  '

# Row 21: Uniswap
$ws.Range("D21").Value = '''6.08'
$ws.Range("E21").Value = '  -4.08%  '

# Row 22: Dai
$ws.Range("E22").Value = '  -0.01%  '

# Row 23: Litecoin
$ws.Range("D23").Value = '''66.32'
$ws.Range("E23").Value = '  -0.45%  '

# Row 24: Polygon
$ws.Range("D24").Value = '''0.424'
$ws.Range("E24").Value = '  +1.79%  '

# Row 25: Binance-PegBSC-USD
$ws.Range("E25").Value = '  -0.08%  '

# Row 26: Kaspa
$ws.Range("E26").Value = '  -4.92%  '

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").Value = '''7.09'
$ws.Range("E27").Value = '  -3.19%  '

# Row 28: PEPE
$ws.Range("D28").Value = '0.0₃0738'
$ws.Range("E28").Value = '  -2.10%  '

# Row 29: USDe
$ws.Range("E29").Value = '  -0.05%  '

# Row 30: PancakeSwap
$ws.Range("D30").Value = '''1.65'
$ws.Range("E30").Value = '  -0.83%  '

# Row 31: Aptos
$ws.Range("E31").Value = '  +1.81%  '

# Row 32: Monero
$ws.Range("D32").Value = '''152.77'
$ws.Range("E32").Value = '  +1.63%  '

# Row 33: EthereumClassic
$ws.Range("D33").Value = '''18.90'
$ws.Range("E33").Value = '  -0.69%  '

# Row 34: NEARProtocol
$ws.Range("D34").Value = '''3.90'
$ws.Range("E34").Value = '  -3.43%  '

# Row 35: SuiNetwork
$ws.Range("D35").Value = '''0.849'
$ws.Range("E35").Value = '  +2.29%  '

# Row 36: ImmutableX
$ws.Range("E36").Value = '  -4.58%  '

# Row 37: Fetch.AI
$ws.Range("D37").Value = '''0.819'
$ws.Range("E37").Value = '  -2.87%  '

# Row 38: Stacks
$ws.Range("E38").Value = '  -3.33%  '

# Row 39: Filecoin
$ws.Range("E39").Value = '  -0.89%  '

# Row 40: Bittensor
$ws.Range("D40").Value = '''278.34'
$ws.Range("E40").Value = '  -5.04%  '

# Row 41: FirstDigitalUSD
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  -0.06%  '

# Row 42: Mantle
$ws.Range("D42").Value = '''0.593'
$ws.Range("E42").Value = '  -2.50%  '

# Row 43: WhiteBITCoin
$ws.Range("D43").Value = '''10.62'
$ws.Range("E43").Value = '  -1.14%  '

# Row 44: Stellar
$ws.Range("D44").Value = '''0.0939'
$ws.Range("E44").Value = '  -1.13%  '

# Row 45: Hedera
$ws.Range("D45").Value = '''0.0527'
$ws.Range("E45").Value = '  -2.89%  '

# Row 46: EnergySwap
$ws.Range("E46").Value = '  -5.55%  '

# Row 47: VeChain
$ws.Range("D47").Value = '''0.0228'
$ws.Range("E47").Value = '  +0.39%  '

# Row 48: Maker
$ws.Range("D48").Value = '1.901.13'
$ws.Range("E48").Value = '  -4.16%  '

# Row 49: InjectiveProtocol
$ws.Range("D49").Value = '''17.85'
$ws.Range("E49").Value = '  -3.33%  '

# Row 50: RenderToken
$ws.Range("E50").Value = '  -3.64%  '

# Row 51: Aave
$ws.Range("D51").Value = '''109.59'
$ws.Range("E51").Value = '  -1.08%  '
